# Trigger GitHub Actions refresh
#
# The nightly PAT-contribution check ran again and appended a new date
# column ("2025-05-11") to the "sample" sheet. Every tracked PAT shows
# "No Contributions" for that day, mirroring the existing C/D columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1: new date header. Format the cell as Text first so Excel stores the
# literal string "2025-05-11" instead of auto-converting it to a date
# serial number (matching how the other date headers - "2025-05-06" /
# "2025-05-10" - are stored as plain text).
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2025-05-11"
$ws.Range("E1").NumberFormat = "General"

# E2:E9: no contributions recorded for any PAT on the new date.
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 5).Value = "No Contributions"
}
